# Replace the "." placeholder text with "NA" for every cell that currently
# holds it. All of these cells share the same underlying shared-string
# entry, so updating all of them collapses back onto a single shared string
# (now reading "NA") instead of leaving the old "." string orphaned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq ".") {
        $cell.Value = "NA"
    }
}

# Restore the recorded UI selection (active cell I8) on the sheet.
$ws.Range("I8").Select()
